# edit.ps1
# Applies data corrections to worksheet "Tab02":
#  1. Fixes mojibake characters in the Regional Economic Communities note (cell A103)
#     (restores proper UTF-8 accented characters: i with acute, i with acute in "Lingua", u with acute in "Comun")
#  2. Updates recalculated/corrected numeric data values for rows 67, 68, 69, 70, 72, 73 (minor precision updates)
#     and rows 97-98 (Africa/ROW Fragile States series), matching the published revision of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab02")

# --- Fix corrupted accented characters in the regional economic communities footnote ---
$ws.Range("A103").Value = 'Regional Economic Communities:CEN-SAD = "Community of Sahel-Saharan States";COMESA = "Common Market for Eastern and Southern Africa";EAC = "East African Community";ECCAS = "Economic Community of Central African States";ECOWAS = "Economic Community of West African States";IGAD = "Intergovernmental Authority on Development";SADC = "Southern African Development Community";UMA = "Arab Maghreb Union";PALOP = "Países Africanos de Língua Oficial Portuguesa";ASEAN = "Association of Southeast Asian Nations";MERCOSUR = "Mercado Común del Sur".EU27 = "European Union (27 members)".OECD = "Organisation for Economic Co-operation and Development".'

# --- Corrected numeric cell values ---
$cellUpdates = @{
    "AJ67" = 4.6837926910574499
    "C67" = 2.5628862630534601
    "V67" = 3.0655261207358402
    "AJ68" = 3.9924538457673799
    "P68" = 6.4957053632020996
    "Q69" = 6.1801702749040199
    "AM70" = 4.7641037426705299
    "O70" = 8.3107341933470806
    "L72" = 4.3615133403029898
    "U72" = 5.1479046481469704
    "AJ73" = 2.29544237607892
    "AA97" = 4.5851264630995301
    "AB97" = 3.81281188181786
    "AC97" = 1.5409000938411399
    "AD97" = 5.0818885118292103
    "AE97" = 3.8292843226646598
    "AF97" = 3.1516989304810101
    "AG97" = -1.3959373897413001
    "AH97" = 5.7085074124332902
    "AI97" = 3.6103000928228699
    "AJ97" = 4.9752604085909899
    "AK97" = 4.8403467960898903
    "AL97" = 4.8493844565952902
    "AM97" = 4.7198496341695604
    "AN97" = 4.9222552284608998
    "AO97" = 4.9198942048660204
    "AP97" = 4.8503201721378497
    "C97" = 2.4110916607942401
    "D97" = 4.0492229204279804
    "E97" = 0.11617974231199001
    "F97" = 0.91884069570705995
    "G97" = 1.82723858306303
    "H97" = 0.95668196446331
    "I97" = 4.9134286107402803
    "J97" = 2.65142820981547
    "K97" = 2.5872635323362001
    "L97" = 2.12746902067525
    "M97" = 3.44401483617705
    "N97" = 5.0723449490637904
    "O97" = 6.5303665369706101
    "P97" = 6.9952529233921501
    "Q97" = 7.4281663813684
    "R97" = 7.3859195267578599
    "S97" = 5.8020410860950999
    "T97" = 7.1609567660121503
    "U97" = 5.6612131878272196
    "V97" = 4.8177129894068598
    "W97" = 8.2001394070253095
    "X97" = 2.1832314580794598
    "Y97" = 12.406156864222
    "Z97" = 4.0575421739085504
    "AA98" = 3.9580308535606199
    "AB98" = 1.8856631093474501
    "AC98" = 4.0017624393951401
    "AD98" = 2.7437632134290402
    "AE98" = 1.9438351990659399
    "AF98" = 1.03327321713176
    "AG98" = 0.35148641209803
    "AH98" = 4.4837751750688399
    "AI98" = 4.8789750798261302
    "AJ98" = 2.8160914204002498
    "AK98" = 3.89964088581692
    "AL98" = 4.3166794501563599
    "AM98" = 4.5236834772465597
    "AN98" = 4.6587809131030102
    "AO98" = 4.6089564575477002
    "AP98" = 4.4011805659257801
    "C98" = 8.7325268008393504
    "D98" = 8.9713894320230292
    "E98" = 5.0616700252218001
    "F98" = 0.49235847999964
    "G98" = 0.34702596254647999
    "H98" = 3.5543719162138898
    "I98" = 4.77453759504077
    "J98" = 2.7442617465422501
    "K98" = 2.6520660614302698
    "L98" = 1.9740126207022199
    "M98" = 5.1671679066532796
    "N98" = 3.3384710749257001
    "O98" = 3.9865498044879502
    "P98" = 5.2912395607747698
    "Q98" = 7.5702920131534404
    "R98" = 6.5568483356751903
    "S98" = 6.5769251305763801
    "T98" = 7.6727495222031896
    "U98" = 3.60824451417967
    "V98" = 1.0407356440943001
    "W98" = 4.2376540312198703
    "X98" = 4.1100317700702798
    "Y98" = 2.1979109512474002
    "Z98" = 2.2649287130377802
}

foreach ($cellRef in $cellUpdates.Keys) {
    $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
}
